# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 5506
$ws1.Range("F8").Value = 905
$ws1.Range("F10").Value = 2453
$ws1.Range("F12").Value = 76
$ws1.Range("F14").Value = 2299
$ws1.Range("F15").Value = 208

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 5506
$ws4.Range("F10").Value = 905
$ws4.Range("F12").Value = 2453
$ws4.Range("F14").Value = 76
$ws4.Range("F17").Value = 2299
$ws4.Range("F18").Value = 208
